# Update Betfair back/lay odds for 2025-12-31 fixtures (Sheet1).
# Only numeric odds cells change; all other cells/rows are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3.75
$ws.Range("G2").Value = 3.85
$ws.Range("H2").Value = 2.16
$ws.Range("I2").Value = 2.18
$ws.Range("J2").Value = 3.55
$ws.Range("N2").Value = 3.7
$ws.Range("O2").Value = 1.34
$ws.Range("P2").Value = 1.93
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 1.35
$ws.Range("S2").Value = 3.6
$ws.Range("T2").Value = 1.81
$ws.Range("U2").Value = 2.12
$ws.Range("V2").Value = 1.84
$ws.Range("X2").Value = 14
$ws.Range("Y2").Value = 9.6
$ws.Range("AA2").Value = 27
$ws.Range("AB2").Value = 14
$ws.Range("AC2").Value = 8.4
$ws.Range("AF2").Value = 28
$ws.Range("AG2").Value = 15.5
$ws.Range("AH2").Value = 18
$ws.Range("AM2").Value = 120
$ws.Range("AO2").Value = 18
$ws.Range("F4").Value = 4.4
$ws.Range("G4").Value = 5.5
$ws.Range("H4").Value = 1.73
$ws.Range("J4").Value = 3.6
$ws.Range("N4").Value = 4.6
$ws.Range("P4").Value = 2.24
$ws.Range("R4").Value = 1.5
$ws.Range("S4").Value = 2.56
$ws.Range("U4").Value = 2.24
$ws.Range("Y4").Value = 990
$ws.Range("AC4").Value = 11.5
$ws.Range("F5").Value = 3.85
$ws.Range("I5").Value = 2.14
$ws.Range("K5").Value = 3.9
$ws.Range("U5").Value = 2.02
$ws.Range("V5").Value = 1.87
$ws.Range("W5").Value = 1.28
$ws.Range("J6").Value = 3.7
$ws.Range("F7").Value = 6.8
$ws.Range("G7").Value = 10.5
$ws.Range("I7").Value = 1.43
$ws.Range("J7").Value = 4.7
$ws.Range("K7").Value = 5.5
$ws.Range("N7").Value = 4.2
$ws.Range("V7").Value = 3
$ws.Range("G8").Value = 3.5
$ws.Range("H8").Value = 2.06
$ws.Range("I8").Value = 2.22
$ws.Range("T8").Value = 1.51
$ws.Range("V8").Value = 1.81
$ws.Range("Y8").Value = 18.5
$ws.Range("Z8").Value = 21
$ws.Range("AA8").Value = 32
$ws.Range("AF8").Value = 34
$ws.Range("AG8").Value = 18
$ws.Range("AI8").Value = 30
$ws.Range("AO8").Value = 11
$ws.Range("F9").Value = 1.97
$ws.Range("G9").Value = 2.08
$ws.Range("H9").Value = 3.6
$ws.Range("I9").Value = 4.1
$ws.Range("J9").Value = 3.95
$ws.Range("N9").Value = 5.4
$ws.Range("P9").Value = 2.48
$ws.Range("R9").Value = 1.59
$ws.Range("S9").Value = 2.38
$ws.Range("T9").Value = 1.55
$ws.Range("U9").Value = 2.46
$ws.Range("V9").Value = 1.33
$ws.Range("W9").Value = 1.92
$ws.Range("Y9").Value = 23
$ws.Range("AA9").Value = 95
$ws.Range("AB9").Value = 15
$ws.Range("AC9").Value = 10.5
$ws.Range("AF9").Value = 17
$ws.Range("AG9").Value = 12
$ws.Range("AH9").Value = 16
$ws.Range("AI9").Value = 980
$ws.Range("AJ9").Value = 28
$ws.Range("AL9").Value = 980
$ws.Range("AN9").Value = 11.5
$ws.Range("AO9").Value = 30
$ws.Range("I10").Value = 1.72
$ws.Range("K10").Value = 4.5
$ws.Range("P10").Value = 2.06
$ws.Range("V10").Value = 2.38
$ws.Range("Y10").Value = 10
$ws.Range("F11").Value = 12
$ws.Range("G11").Value = 15
$ws.Range("K11").Value = 8
$ws.Range("Q11").Value = 1.35
$ws.Range("T11").Value = 1.78
$ws.Range("W11").Value = 1.07
$ws.Range("AJ11").Value = 480
$ws.Range("AK11").Value = 190
$ws.Range("F12").Value = 4.7
$ws.Range("H12").Value = 1.68
$ws.Range("I12").Value = 1.79
$ws.Range("J12").Value = 4.1
$ws.Range("K12").Value = 4.7
$ws.Range("Q12").Value = 1.63
$ws.Range("T12").Value = 1.67
$ws.Range("V12").Value = 2.24
$ws.Range("F13").Value = 1.28
$ws.Range("H13").Value = 11
$ws.Range("I13").Value = 13.5
$ws.Range("J13").Value = 6
$ws.Range("K13").Value = 7.2
$ws.Range("L13").Value = 1.19
$ws.Range("O13").Value = 1.17
$ws.Range("Q13").Value = 1.54
$ws.Range("R13").Value = 1.64
$ws.Range("S13").Value = 2.26
$ws.Range("V13").Value = 1.08
$ws.Range("W13").Value = 3.95
$ws.Range("AN13").Value = 4.4
